$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated values for columns H (TNA BANCARIA), J (Interes del Banco por el
# plazo elegido), K (Precio final financiado) and L (TNA con todos los
# gastos incluidos) for the affected rows.
$updates = @{
    41  = @{ H = 0.47099999999999997; J = 23.227397260273975; K = 129.7584493150685;  L = 0.60343522222222234 }
    42  = @{ H = 0.53520000000000001; J = 35.191232876712327; K = 142.35636821917808; L = 0.64416976666666659 }
    43  = @{ H = 0.56520000000000004; J = 27.872876712328765; K = 134.65013917808218; L = 0.70262782222222231 }
    44  = @{ H = 0.56520000000000004; J = 37.163835616438355; K = 144.4335189041096;  L = 0.67575976666666682 }
    48  = @{ H = 0.53;                J = 26.136986301369863; K = 126.76767123287669; L = 0.54278888888888843 }
    49  = @{ H = 0.53;                J = 39.205479452054796; K = 139.90150684931504; L = 0.53940925925925887 }
    50  = @{ H = 0.53;                J = 52.273972602739725; K = 153.03534246575342; L = 0.53771944444444431 }
    51  = @{ H = 0.42;                J = 20.712328767123285; K = 127.06560922855081; L = 0.54883040935672478 }
    52  = @{ H = 0.44500000000000001; J = 32.917808219178077; K = 139.91348233597691; L = 0.53957115009746559 }
    53  = @{ H = 0.46500000000000002; J = 45.863013698630141; K = 153.54001441961066; L = 0.54283625730994123 }
    54  = @{ H = 0.28999999999999998; J = 14.301369863013697; K = 127.00152207001521; L = 0.54753086419753061 }
    55  = @{ H = 0.35499999999999998; J = 26.260273972602739; K = 140.28919330289193; L = 0.54465020576131684 }
    56  = @{ H = 0.38500000000000001; J = 37.972602739726028; K = 153.30289193302892; L = 0.54043209876543197 }
    57  = @{ H = 0.28999999999999998; J = 14.301369863013697; K = 126.7602191780822;  L = 0.542637777777778   }
    58  = @{ H = 0.36499999999999999; J = 27;                 K = 140.84300000000002; L = 0.55213685185185202 }
    59  = @{ H = 0.4;                 J = 39.452054794520549; K = 154.65232876712329; L = 0.55411388888888891 }
    60  = @{ H = 0.32;                J = 15.780821917808218; K = 128.4009315068493;  L = 0.57590777777777757 }
    61  = @{ H = 0.39;                J = 28.849315068493155; K = 142.8938904109589;  L = 0.5798618518518518  }
    62  = @{ H = 0.42499999999999999; J = 41.917808219178085; K = 157.3868493150685;  L = 0.58183888888888902 }
    116 = @{ H = 0.54;                J = 8.8767123287671232; K = 111.05424657534246; L = 0.67246666666666655 }
    117 = @{ H = 0.54;                J = 13.315068493150687; K = 115.58136986301371; L = 0.63191111111111142 }
    118 = @{ H = 0.54;                J = 26.630136986301373; K = 129.16273972602738; L = 0.59135555555555541 }
    119 = @{ H = 0.54;                J = 39.945205479452063; K = 142.7441095890411;  L = 0.57783703703703726 }
    120 = @{ H = 0.54;                J = 53.260273972602747; K = 156.32547945205479; L = 0.5710777777777778  }
    126 = @{ H = 0.47;                J = 7.7260273972602738; K = 109.88054794520548; L = 0.60106666666666642 }
    127 = @{ H = 0.47;                J = 11.58904109589041;  K = 113.82082191780822; L = 0.56051111111111129 }
    128 = @{ H = 0.47;                J = 23.17808219178082;  K = 125.64164383561643; L = 0.51995555555555573 }
    129 = @{ H = 0.47;                J = 34.767123287671232; K = 137.46246575342465; L = 0.50643703703703691 }
    130 = @{ H = 0.47;                J = 46.356164383561641; K = 149.28328767123287; L = 0.49967777777777767 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 8).Value  = $vals.H   # Column H
    $ws.Cells.Item($row, 10).Value = $vals.J   # Column J
    $ws.Cells.Item($row, 11).Value = $vals.K   # Column K
    $ws.Cells.Item($row, 12).Value = $vals.L   # Column L
}
